$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in A2 from "tiger" to "ansul"
$ws.Range("A2").Value = "ansul"

# Update the active selection cell from C8 to D8
$ws.Range("D8").Select()
